$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R9").Style = "Normal"

$ws.Cells.Item(9, 2).Value = "id"
for ($i = 1; $i -le 141; $i++) {
    $ws.Cells.Item(9, $i + 2).Value = $i
}
